$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (cell, new value). Values for the Price (D) and Volume(1h) (E) columns
# are prefixed with a literal apostrophe so Excel stores them as text -- matching the
# original inlineStr cells -- instead of auto-converting "317.94" to a number or
# "3.79%" to 0.0379.
$edits = @(
    ,("D2", '''317.94')
    ,("E2", '''3.79%')
    ,("D3", '''39.77')
    ,("E3", '''2.21%')
    ,("D4", '''5.143')
    ,("E4", '''1.01%')
    ,("D5", '''0.08200')
    ,("E5", '''1.79%')
    ,("D6", '''2.114')
    ,("E6", '''8.95%')
    ,("D7", '''8.304')
    ,("E7", '''3.81%')
    ,("B8", 'GateToken')
    ,("C8", 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt')
    ,("D8", '''4.289')
    ,("E8", '''2.27%')
    ,("B9", 'MXToken')
    ,("C9", 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx')
    ,("D9", '''0.9350')
    ,("E9", '''0.36%')
    ,("B10", 'LiechtensteinCryptoassetsExchange')
    ,("C10", 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx')
    ,("D10", '''0.1377')
    ,("E10", '''-5.45%')
    ,("B11", 'WazirX')
    ,("C11", 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx')
    ,("D11", '''0.1987')
    ,("E11", '''3.11%')
    ,("B12", 'MandalaExchangeToken')
    ,("C12", 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx')
    ,("D12", '''0.09104')
    ,("E12", '''0.97%')
    ,("B13", 'BitrueCoin')
    ,("C13", 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr')
    ,("D13", '''0.03480')
    ,("E13", '''-0.59%')
    ,("B14", 'BitMartToken')
    ,("C14", 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx')
    ,("D14", '''0.09799')
    ,("E14", '''0.18%')
    ,("B15", 'BitForexToken')
    ,("C15", 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf')
    ,("D15", '''0.001399')
    ,("E15", '''-0.45%')
    ,("B16", 'TigerCash')
    ,("C16", 'https://coinranking.com/coin/6hIn06L2+tigercash-tch')
    ,("D16", '''0.006077')
    ,("E16", '''4.49%')
    ,("B17", 'LEO')
    ,("C17", 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo')
    ,("D17", '''3.679')
    ,("E17", '''-2.93%')
    ,("E18", '''-3.67%')
    ,("D19", '''0.3473')
    ,("E19", '''1.56%')
    ,("E20", '''-3.10%')
    ,("D21", '''4.901')
    ,("E21", '''2.35%')
    ,("D22", '''0.2449')
    ,("E22", '''1.49%')
    ,("D23", '''0.04325')
    ,("E23", '''-1.31%')
    ,("D24", '''0.001226')
    ,("E24", '''-0.98%')
    ,("D25", '''0.004770')
    ,("E25", '''11.46%')
    ,("D26", '''0.0001300')
    ,("E26", '''-0.11%')
    ,("D27", '''0.0003998')
    ,("E27", '''-10.12%')
    ,("D39", '''0.02234')
    ,("E39", '''9.77%')
    ,("D40", '''0.05225')
    ,("E40", '''3.85%')
    ,("D41", '''0.007553')
    ,("E41", '''1.43%')
    ,("D42", '''0.009671')
    ,("E42", '''-6.15%')
    ,("D43", '''0.1384')
    ,("E43", '''2.67%')
    ,("D44", '''0.002150')
    ,("E44", '''1.30%')
    ,("D45", '''0.009185')
    ,("E45", '''1.63%')
    ,("D46", '''0.00006610')
    ,("E46", '''6.49%')
    ,("D47", '''0.00000000750')
    ,("E47", '''-0.11%')
    ,("B48", 'CoinbaseStockToken')
    ,("C48", 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin')
    ,("D48", '''0.001200')
    ,("E48", '''-25.07%')
    ,("B49", 'BOLO')
    ,("C49", 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo')
    ,("D49", '''0.002775')
    ,("E49", '''0.22%')
    ,("D50", '''0.00002100')
    ,("E50", '''-0.11%')
    ,("D51", '''0.0002000')
    ,("E51", '''-0.11%')
)

foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value = $edit[1]
}
